# Estado de Cuenta - add new dues period (2509) data, update summary totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 6 new blank rows right after the current last data row (93),
#    before the footer (signature) block. This shifts the footer rows
#    (98:99) down to (104:105) and keeps merged cells intact.
# ---------------------------------------------------------------------------

# Remember the "closing" (bottom-of-table) formatting of the current last
# row (93) so it can be re-applied to the new last row once the table grows.
$ws.Range("B93:J93").Copy()
$closingFormatRow = $ws.Range("B200:J200")
$closingFormatRow.PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B94:B99").EntireRow.Insert(-4121)  # xlShiftDown

# Copy the regular (interior) row formatting down onto the six new rows,
# plus row 93 (which is no longer the closing row of the table) so every
# row from 93 to 98 matches the rest of the table (borders, currency
# format, etc.)
$ws.Range("B92:J92").Copy()
$ws.Range("B93:J98").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = $false

# Row 99 is now the new closing (last) row of the table - give it back the
# special "bottom of table" formatting that row 93 used to have.
$closingFormatRow.Copy()
$ws.Range("B99:J99").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = $false
$closingFormatRow.EntireRow.ClearFormats()

# ---------------------------------------------------------------------------
# 2) Fill in the six new rows for period 2509.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Doc="1047377858";  Name="EFRANIO MARTELO BELTRAN";              Mora=67759; Salario=1693965 },
    @{ Doc="73352492";    Name="DAVID JOSE DIAZ CANTILLO";             Mora=56940; Salario=1423500 },
    @{ Doc="1066184330";  Name="JUAN DANIEL PATERNINA MENDOZA";        Mora=67759; Salario=1693965 },
    @{ Doc="1007521765";  Name="ANYIE PAHOLA PACHECO OSORIO";          Mora=56940; Salario=1423500 },
    @{ Doc="8866266";     Name="ARTURO ELIAS TORRES CORREA";           Mora=44000; Salario=1100000 },
    @{ Doc="1143392420";  Name="LAURA GENOVEVA MENDOZA CUADRO";        Mora=56940; Salario=1423500 }
)

$r = 94
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $row.Doc
    $ws.Cells.Item($r, 4).Value = $row.Name
    $ws.Cells.Item($r, 5).Value = "2509"
    $ws.Cells.Item($r, 6).Value = $row.Mora
    $ws.Cells.Item($r, 7).Value = $row.Salario
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Update the summary header figures.
#    - Valor Mora total (sum of all "Valor Mora" rows)
#    - Cant. Trabajadores (unique workers)
#    - Cant. Periodos (unique periods)
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 5).Value = 3809746
$ws.Cells.Item(13, 3).Value = 16
$ws.Cells.Item(13, 6).Value = 42
